$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 181
$ws.Range("A2").Value = 200
$ws.Range("A3").Value = 200
$ws.Range("A4").Value = 189
$ws.Range("A5").Value = 177
$ws.Range("A6").Value = 187
$ws.Range("A7").Value = 203.3999999999996
$ws.Range("A8").Value = 200
